$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update algorithm name/version headers (D1, E1, F1)
$ws.Range("D1").Value = "Dilithium2"
$ws.Range("E1").Value = "Falcon-512"
$ws.Range("F1").Value = "SPHINCS+-SHA2-128f-simple"

# Row 4 - signatureSize
$ws.Range("B4").Value = 128.0
$ws.Range("C4").Value = 72.0
$ws.Range("D4").Value = 2420.0
$ws.Range("E4").Value = 654.0
$ws.Range("F4").Value = 17088.0

# Row 5 - publicKeySize
$ws.Range("B5").Value = 162.0
$ws.Range("C5").Value = 91.0
$ws.Range("D5").Value = 1312.0
$ws.Range("E5").Value = 897.0
$ws.Range("F5").Value = 32.0

# Row 6 - privateKeySize
$ws.Range("B6").Value = 635.0
$ws.Range("C6").Value = 150.0
$ws.Range("D6").Value = 2528.0
$ws.Range("E6").Value = 1281.0
$ws.Range("F6").Value = 64.0

# Row 7 - keyGenTimeMean
$ws.Range("B7").Value = 148186112.0
$ws.Range("C7").Value = 2178852.0
$ws.Range("D7").Value = 35983.0
$ws.Range("E7").Value = 9695598.0
$ws.Range("F7").Value = 333133.0

# Row 9 - signatureTimeMean
$ws.Range("B9").Value = 648755.0
$ws.Range("C9").Value = 2924291.0
$ws.Range("D9").Value = 50663.0
$ws.Range("E9").Value = 446189.0
$ws.Range("F9").Value = 8258050.0

# Row 11 - verifyTimeMean
$ws.Range("B11").Value = 88476.0
$ws.Range("C11").Value = 4067236.0
$ws.Range("D11").Value = 34141.0
$ws.Range("E11").Value = 55766.0
$ws.Range("F11").Value = 847349.0
